$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Design")

# Print i2c devices while setup:
# The 8x20 grid (rows 17-20, cols C:V) is rewritten so each row now
# repeats two sequential values (1,2 / 3,4 / 5,6 / 7,8) across the row
# instead of the previous repeating 1-8 cycle.
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

for ($r = 17; $r -le 20; $r++) {
    $base = ($r - 17) * 2
    for ($i = 0; $i -lt $cols.Length; $i++) {
        if ($i % 2 -eq 0) {
            $val = $base + 1
        } else {
            $val = $base + 2
        }
        $ws.Range($cols[$i] + $r).Value = $val
    }
}

# Move the sheet's active selection from X2:Z2 to S8
[void]$ws.Range("S8").Select()
